$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.288.47'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.157.33'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.88'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.12'
$ws.Range("E6").Value = '  -4.15%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -4.81%  '
$ws.Range("E9").Value = '  -3.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.713.00'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '64.395.60'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.28'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.155.60'
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '405.99'
$ws.Range("E18").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.27'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.09'
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.45'
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.198'
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("E26").Value = '  -7.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.85'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.991'
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("E30").Value = '  -3.49%  '
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.98'
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.686.83'
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.91'
$ws.Range("E38").Value = '  -5.35%  '
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.695'
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.40'
$ws.Range("E42").Value = '  -6.21%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0258'
$ws.Range("E43").Value = '  -1.88%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '291.27'
$ws.Range("E44").Value = '  -2.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.33'
$ws.Range("E45").Value = '  -4.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0983'
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("E48").Value = '  -10.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.47'
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.881'
$ws.Range("E51").Value = '  -5.34%  '
